# Apply targeted updates to the weekly work report workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# Update the "Report Generated On" timestamp in D5.
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

# Zero out the billing/pricing figures (no billable work this period).
$ws.Range("C8").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
